# Updated symbol list on Fri Feb 17 05:27:25 UTC 2023 with GitHub Actions
# Applies updated crypto price/volume data and swaps rows 11 and 12
# (MandalaExchangeToken <-> BitrueCoin) to match the refreshed source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.39"
$ws.Range("E2").Value = "'-5.21%"
$ws.Range("D3").Value = "'48.58"
$ws.Range("E3").Value = "'-2.32%"
$ws.Range("D4").Value = "'5.162"
$ws.Range("E4").Value = "'-3.27%"
$ws.Range("D5").Value = "'0.07757"
$ws.Range("E5").Value = "'-4.95%"
$ws.Range("D6").Value = "'4.507"
$ws.Range("E6").Value = "'-2.31%"
$ws.Range("D7").Value = "'1.341"
$ws.Range("E7").Value = "'14.34%"
$ws.Range("D8").Value = "'1.552"
$ws.Range("E8").Value = "'-6.68%"
$ws.Range("E9").Value = "'-9.39%"
$ws.Range("D10").Value = "'0.1929"
$ws.Range("E10").Value = "'-1.30%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.04653"
$ws.Range("E11").Value = "'2.30%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09277"
$ws.Range("E12").Value = "'-2.68%"
$ws.Range("D13").Value = "'0.1048"
$ws.Range("E13").Value = "'-0.03%"
$ws.Range("D14").Value = "'0.001264"
$ws.Range("E14").Value = "'-5.20%"
$ws.Range("D15").Value = "'0.04184"
$ws.Range("E15").Value = "'-2.71%"
$ws.Range("D16").Value = "'0.005872"
$ws.Range("E16").Value = "'-1.20%"
$ws.Range("E17").Value = "'-1.94%"
$ws.Range("D18").Value = "'2.276"
$ws.Range("E18").Value = "'-6.67%"
$ws.Range("E19").Value = "'2.76%"
$ws.Range("D20").Value = "'8.043"
$ws.Range("E20").Value = "'-1.71%"
$ws.Range("D21").Value = "'0.1340"
$ws.Range("E21").Value = "'-5.61%"
$ws.Range("D22").Value = "'0.3039"
$ws.Range("E22").Value = "'-0.46%"
$ws.Range("E23").Value = "'-2.47%"
$ws.Range("D24").Value = "'0.004085"
$ws.Range("E24").Value = "'-4.18%"
$ws.Range("E25").Value = "'0.19%"
$ws.Range("E26").Value = "'-4.06%"
$ws.Range("D38").Value = "'0.02557"
$ws.Range("E38").Value = "'-8.13%"
$ws.Range("D39").Value = "'0.05827"
$ws.Range("E39").Value = "'4.81%"
$ws.Range("D40").Value = "'0.01076"
$ws.Range("E40").Value = "'70.77%"
$ws.Range("D41").Value = "'0.007924"
$ws.Range("E41").Value = "'3.14%"
$ws.Range("D42").Value = "'0.1420"
$ws.Range("E42").Value = "'-1.99%"
$ws.Range("D43").Value = "'0.008347"
$ws.Range("E43").Value = "'8.45%"
$ws.Range("D44").Value = "'0.007714"
$ws.Range("E44").Value = "'-4.51%"
$ws.Range("D45").Value = "'0.3069"
$ws.Range("E45").Value = "'-12.81%"
$ws.Range("D46").Value = "'0.00006987"
$ws.Range("E46").Value = "'3.11%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.21%"
$ws.Range("E48").Value = "'-7.54%"
$ws.Range("E49").Value = "'0.13%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.21%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.21%"
